# Auto update stock data
# Update the "current" date (first row of each 6-row stock block) from
# 2025/12/07 to 2025/12/08 in column A, keeping the cell as plain text
# (not auto-converted into a date serial number by Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "2025/12/07") {
        # Force the value to be stored as text (not parsed as a date)
        $cell.NumberFormat = "@"
        $cell.Value2 = "2025/12/08"
        # Restore the default/normal cell style so no extra formatting is introduced
        $cell.Style = "Normal"
    }
}
